# Update the month_task service schedule: shift from December to August
# (new date ranges, reassigned servers, updated scripture readings).
$d = $word.ActiveDocument
$failures = 0

$found = $d.Content.Find.Execute("12月主日崇拜服侍表", $true, $false, $false, $false, $false, $true, 1, $false, "8月主日崇拜服侍表", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #0" }
$found = $d.Content.Find.Execute("日期         2日           9日           16日          23日          30日          ", $true, $false, $false, $false, $false, $true, 1, $false, "日期         5日           12日          19日          26日          ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #1" }
$found = $d.Content.Find.Execute("领诗         昀晏         若涵         君昊         卢伟         思倩         ", $true, $false, $false, $false, $false, $true, 1, $false, "领诗         若涵         柏万         李志         君昊         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #2" }
$found = $d.Content.Find.Execute("司乐         君昊         泰禾         昀晏         曜宏         靖珊         ", $true, $false, $false, $false, $false, $true, 1, $false, "司乐         柏万         昀晏         泰禾         曜宏         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #3" }
$found = $d.Content.Find.Execute("音控         成禾         君昊         燕杉         明亮         卢伟         ", $true, $false, $false, $false, $false, $true, 1, $false, "音控         以恒         灿荣         君昊         成禾         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #4" }
$found = $d.Content.Find.Execute("司会         元鹏         兴伟         卢伟         郑凯         明亮         ", $true, $false, $false, $false, $false, $true, 1, $false, "司会         元鹏         李志         曜宏         兴伟         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #5" }
$found = $d.Content.Find.Execute("圣餐         明亮         无            无            无            无            ", $true, $false, $false, $false, $false, $true, 1, $false, "圣餐         卢伟         无            无            无            ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #6" }
$found = $d.Content.Find.Execute("讲道         蔡牧师      蔡牧师      蔡牧师      蔡牧师      蔡牧师      ", $true, $false, $false, $false, $false, $true, 1, $false, "讲道         蔡牧师      蔡牧师      蔡牧师      蔡牧师      ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #7" }
$found = $d.Content.Find.Execute("茶点         李志         成禾         卓阳         佳宁         刘欢         ", $true, $false, $false, $false, $false, $true, 1, $false, "茶点         志钰         自超         董婷         燕杉         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #8" }
$found = $d.Content.Find.Execute("打扫         李志         成禾         卓阳         佳宁         刘欢         ", $true, $false, $false, $false, $false, $true, 1, $false, "打扫         灿荣         自超         董婷         燕杉         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #9" }
$found = $d.Content.Find.Execute("接待         若涵         佳宁         李志         兴伟         自超         ", $true, $false, $false, $false, $false, $true, 1, $false, "接待         董婷         君昊         灿荣         佳宁         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #10" }
$found = $d.Content.Find.Execute("儿童         蔡师母      文佳         悦茜         蔡师母      文佳         ", $true, $false, $false, $false, $false, $true, 1, $false, "儿童         蔡师母      文佳         岩美         悦茜         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #11" }
$found = $d.Content.Find.Execute("助手         兴伟         以恒         曜宏         天洵         兴伟         ", $true, $false, $false, $false, $false, $true, 1, $false, "助手         天洵         蔡师母      宗尉         蔡师母      ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #12" }
$found = $d.Content.Find.Execute("12月周五查经服侍表", $true, $false, $false, $false, $false, $true, 1, $false, "8月周五查经服侍表", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #13" }
$found = $d.Content.Find.Execute("日期         7日           14日          21日          28日          ", $true, $false, $false, $false, $false, $true, 1, $false, "日期         3日           10日          17日          24日          31日          ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #14" }
$found = $d.Content.Find.Execute("领诗         以恒         李志         若涵         卓阳         ", $true, $false, $false, $false, $false, $true, 1, $false, "领诗         以恒         卢伟         董婷         蔡师母      曜宏         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #15" }
$found = $d.Content.Find.Execute("司乐         蔡师母      曜宏         泰禾         蔡师母      ", $true, $false, $false, $false, $false, $true, 1, $false, "司乐         泰禾         柏万         蔡师母      曜宏         泰禾         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #16" }
$found = $d.Content.Find.Execute("带领1        李志         佳宁         董婷         悦茜         ", $true, $false, $false, $false, $false, $true, 1, $false, "带领1        曜宏         灿荣         兴伟         君昊         明星         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #17" }
$found = $d.Content.Find.Execute("带领2        卢伟         兴伟         明亮         曜宏         ", $true, $false, $false, $false, $false, $true, 1, $false, "带领2        悦茜         明亮         佳宁         卢伟         李志         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #18" }
$found = $d.Content.Find.Execute("经文         创50章       TBD            TBD            TBD            ", $true, $false, $false, $false, $false, $true, 1, $false, "经文         创32章       创33章       创34章       创35章       创36章       ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #19" }
$found = $d.Content.Find.Execute("茶点         成禾         自超         刘欢         卓阳         ", $true, $false, $false, $false, $false, $true, 1, $false, "茶点         刘欢         李志         蔡师母      佳宁         成禾         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #20" }
$found = $d.Content.Find.Execute("打扫         成禾         自超         刘欢         卓阳         ", $true, $false, $false, $false, $false, $true, 1, $false, "打扫         刘欢         李志         蔡师母      佳宁         成禾         ", 2)
if (-not $found) { $failures++; Write-Output "MISSING MATCH #21" }

Write-Output "Replacements applied with $failures failure(s)."
